# Refresh the coinranking Price (D) / Volume(1h) (E) columns, and for the
# four rows whose ranking changed, the Coin name (B) and Link (C) too -
# matching the GitHub Actions "Updated cryptos list" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '19.956.92'
$ws.Range("E2").Value = '  -8.20%  '
$ws.Range("D3").Value = '1.414.32'
$ws.Range("E3").Value = '  -8.13%  '
$ws.Range("D4").Value = '''1.002'
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("E5").Value = '  +0.13%  '
$ws.Range("D6").Value = '''272.86'
$ws.Range("E6").Value = '  -5.88%  '
$ws.Range("D7").Value = '''0.3690'
$ws.Range("E7").Value = '  -5.19%  '
$ws.Range("D8").Value = '''0.3064'
$ws.Range("E8").Value = '  -3.85%  '
$ws.Range("D9").Value = '''39.16'
$ws.Range("E9").Value = '  -8.76%  '
$ws.Range("E10").Value = '  -5.53%  '
$ws.Range("D11").Value = '''0.06552'
$ws.Range("E11").Value = '  -9.04%  '
$ws.Range("D12").Value = '''1.003'
$ws.Range("E12").Value = '  +0.17%  '
$ws.Range("D13").Value = '''5.419'
$ws.Range("E13").Value = '  -4.02%  '
$ws.Range("D14").Value = '''16.96'
$ws.Range("E14").Value = '  -9.01%  '
$ws.Range("D15").Value = '''6.151'
$ws.Range("E15").Value = '  -7.02%  '
$ws.Range("D16").Value = '1.416.22'
$ws.Range("E16").Value = '  -8.34%  '
$ws.Range("D17").Value = '''0.00001005'
$ws.Range("E17").Value = '  -9.31%  '
$ws.Range("D18").Value = '''0.05739'
$ws.Range("E18").Value = '  -12.88%  '
$ws.Range("D19").Value = '''73.95'
$ws.Range("E19").Value = '  -11.20%  '
$ws.Range("E20").Value = '  +0.14%  '
$ws.Range("D21").Value = '''5.592'
$ws.Range("E21").Value = '  -9.16%  '
$ws.Range("D22").Value = '''14.39'
$ws.Range("E22").Value = '  -6.55%  '
$ws.Range("D23").Value = '''10.84'
$ws.Range("E23").Value = '  -0.33%  '
$ws.Range("D24").Value = '''2.328'
$ws.Range("E24").Value = '  -2.61%  '
$ws.Range("D25").Value = '19.970.90'
$ws.Range("E25").Value = '  -8.16%  '
$ws.Range("D26").Value = '''2.265'
$ws.Range("E26").Value = '  -4.64%  '
$ws.Range("D27").Value = '''139.38'
$ws.Range("E27").Value = '  -4.95%  '
$ws.Range("D28").Value = '''16.91'
$ws.Range("E28").Value = '  -8.03%  '
$ws.Range("D29").Value = '1.578.88'
$ws.Range("E29").Value = '  -8.06%  '
$ws.Range("D30").Value = '''108.58'
$ws.Range("E30").Value = '  -7.65%  '
$ws.Range("E31").Value = '  -19.31%  '
$ws.Range("D32").Value = '''5.348'
$ws.Range("E32").Value = '  -9.86%  '
$ws.Range("D33").Value = '''0.8529'
$ws.Range("E33").Value = '  -12.83%  '
$ws.Range("D34").Value = '''0.07700'
$ws.Range("E34").Value = '  -6.16%  '
$ws.Range("D35").Value = '''8.386'
$ws.Range("E35").Value = '  -5.03%  '
$ws.Range("D36").Value = '''0.05730'
$ws.Range("E36").Value = '  -5.89%  '
$ws.Range("D37").Value = '''1.002'
$ws.Range("E37").Value = '  +0.13%  '
$ws.Range("D38").Value = '''4.766'
$ws.Range("E38").Value = '  -7.30%  '
$ws.Range("D39").Value = '''10.66'
$ws.Range("E39").Value = '  -0.31%  '
$ws.Range("B40").Value = 'Algorand'
$ws.Range("C40").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D40").Value = '''0.1917'
$ws.Range("E40").Value = '  -6.02%  '
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").Value = '''0.02024'
$ws.Range("E41").Value = '  -8.22%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").Value = '''1.062'
$ws.Range("E42").Value = '  -10.77%  '
$ws.Range("B43").Value = 'WEMIXTOKEN'
$ws.Range("C43").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D43").Value = '''1.275'
$ws.Range("E43").Value = '  -13.70%  '
$ws.Range("D44").Value = '''0.5277'
$ws.Range("E44").Value = '  -8.30%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '''12.28'
$ws.Range("E45").Value = '  -5.99%  '
$ws.Range("B46").Value = 'PancakeSwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D46").Value = '''3.528'
$ws.Range("E46").Value = '  -5.81%  '
$ws.Range("D47").Value = '''0.5097'
$ws.Range("E47").Value = '  -7.54%  '
$ws.Range("D48").Value = '''1.797'
$ws.Range("E48").Value = '  -4.03%  '
$ws.Range("D49").Value = '''109.00'
$ws.Range("E49").Value = '  -6.92%  '
$ws.Range("D50").Value = '''1.046'
$ws.Range("E50").Value = '  -10.11%  '
$ws.Range("D51").Value = '''1.003'
$ws.Range("E51").Value = '  +0.19%  '
